# Applies the "nuevos layouts y actualiza backlog" update to the BackLog sheet.
#
# Summary of the change:
#  - Column H ("Comentarios") gets new status notes for several requirement rows.
#  - The "Requerimientos" text for the "Realizar solicitud" story (row 7) is
#    reworded (TextView/EditText wording simplified).
#  - The active selection is moved to H13.
#  - Row 7 grows taller to fit the reworded (slightly longer) text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Comentarios" (column H) notes.
$ws.Cells.Item(2, 8).Value  = "Layout OK"
$ws.Cells.Item(3, 8).Value  = "Pendiente Toast"

# Reword the "Requerimientos" (column G) text for row 7.
$ws.Cells.Item(7, 7).Value  = 'Activity con 1 EditText donde se le solicita que en el  especifique tipo de carga y peso, un boton de solicitar servicio, 1 Alert Dialoge con mensaje que diga "Solicitud recibida su paquete será recogido el día de mañana"'
$ws.Cells.Item(7, 8).Value  = "Pendiente Alert Dialoge"

$ws.Cells.Item(11, 8).Value = "Layout OK"
$ws.Cells.Item(12, 8).Value = "Layout OK"

# Row 7 now needs more vertical room for the updated requirement text.
$ws.Rows.Item(7).RowHeight = 119.4

# Leave the cursor/selection on H13, matching where the author ended up.
$ws.Range("H13").Select()
